$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.269.71"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "2.366.80"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'318.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").Value = "'108.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.91%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").Value = "'0.614"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("D10").Value = "'40.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.20%  "
$ws.Range("D11").Value = "'0.0918"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.78%  "
$ws.Range("D12").Value = "'8.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.02%  "
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").Value = "'0.980"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.61%  "
$ws.Range("D15").Value = "2.727.18"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("E16").Value = "  -4.01%  "
$ws.Range("D17").Value = "2.367.20"
$ws.Range("E17").Value = "  -2.60%  "
$ws.Range("D18").Value = "45.186.07"
$ws.Range("E18").Value = "  -1.40%  "
$ws.Range("D19").Value = "'15.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +15.26%  "
$ws.Range("D20").Value = "'7.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.69%  "
$ws.Range("E21").Value = "  -2.58%  "
$ws.Range("D22").Value = "'3.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").Value = "'73.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("D24").Value = "'265.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("D30").Value = "'22.38"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.27%  "
$ws.Range("E31").Value = "  -2.37%  "
$ws.Range("D32").Value = "'37.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.74%  "
$ws.Range("D33").Value = "'168.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("E34").Value = "  -3.91%  "
$ws.Range("E35").Value = "  -0.57%  "
$ws.Range("E36").Value = "  -4.72%  "
$ws.Range("E37").Value = "  -5.78%  "
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").Value = "'4.05"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.90%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.65%  "
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("D42").Value = "'99.43"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.98%  "
$ws.Range("D43").Value = "'70.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").Value = "1.881.56"
$ws.Range("E44").Value = "  +13.40%  "
$ws.Range("D45").Value = "'12.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.95%  "
$ws.Range("E46").Value = "  -5.32%  "
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").Value = "'5.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").Value = "'84.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.58%  "
$ws.Range("D50").Value = "'111.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.40%  "
$ws.Range("D51").Value = "'9.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.80%  "
